$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 899.68115
$ws.Range("I112").Value = 645
$ws.Range("J112").Value = 923.9365
$ws.Range("K112").Value = 1935
$ws.Range("L112").Value = 2771.8095
$ws.Range("M112").Value = -827
$ws.Range("N112").Value = -4987.8095
$ws.Range("H132").Value = 1025.6316
$ws.Range("I132").Value = 1026.1621
$ws.Range("K132").Value = 3078.4863
$ws.Range("M132").Value = -548.4863
$ws.Range("H137").Value = 1764.0358
$ws.Range("I137").Value = 1209
$ws.Range("J137").Value = 2072.389
$ws.Range("K137").Value = 3627
$ws.Range("L137").Value = 6217.167
$ws.Range("M137").Value = -1077
$ws.Range("N137").Value = -11317.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2688.889
$ws.Range("I2").Value = 3100
$ws.Range("J2").Value = 1250
$ws.Range("K2").Value = 3100
$ws.Range("L2").Value = 1250
$ws.Range("M2").Value = -2987
$ws.Range("N2").Value = -1476
$ws.Range("H32").Value = 3761.8076
$ws.Range("I32").Value = 3114.7231
$ws.Range("J32").Value = 6997.231
$ws.Range("K32").Value = 3114.7231
$ws.Range("L32").Value = 6997.231
$ws.Range("M32").Value = -2827.7231
$ws.Range("N32").Value = -7571.231
$ws.Range("H63").Value = 55557628
$ws.Range("I63").Value = 111112856
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 111112856
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -111112170
$ws.Range("N63").Value = -3772
$ws.Range("H66").Value = 55557628
$ws.Range("I66").Value = 111112856
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 555564280
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -555560848
$ws.Range("N66").Value = -18864
$ws.Range("H116").Value = 2688.889
$ws.Range("I116").Value = 3100
$ws.Range("J116").Value = 1250
$ws.Range("K116").Value = 3100
$ws.Range("L116").Value = 1250
$ws.Range("M116").Value = -806
$ws.Range("N116").Value = -5838

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2688.889
$ws.Range("I3").Value = 3100
$ws.Range("J3").Value = 1250
$ws.Range("K3").Value = 3100
$ws.Range("L3").Value = 1250
$ws.Range("M3").Value = -2986
$ws.Range("N3").Value = -1478
$ws.Range("H99").Value = 111112700
$ws.Range("I99").Value = 142858380
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 142858380
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -142856882
$ws.Range("N99").Value = -5796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 10990477
$ws.Range("I16").Value = 12822060
$ws.Range("J16").Value = 980
$ws.Range("K16").Value = 12822060
$ws.Range("L16").Value = 980
$ws.Range("M16").Value = -12821773
$ws.Range("N16").Value = -1554
$ws.Range("H31").Value = 2200.3765
$ws.Range("I31").Value = 1513.7368
$ws.Range("J31").Value = 2398.0454
$ws.Range("K31").Value = 1513.7368
$ws.Range("L31").Value = 2398.0454
$ws.Range("M31").Value = -1218.7368
$ws.Range("N31").Value = -2988.0454
$ws.Range("H34").Value = 2200.3765
$ws.Range("I34").Value = 1513.7368
$ws.Range("J34").Value = 2398.0454
$ws.Range("K34").Value = 1513.7368
$ws.Range("L34").Value = 2398.0454
$ws.Range("M34").Value = -1311.7368
$ws.Range("N34").Value = -2802.0454
$ws.Range("H86").Value = 71430590
$ws.Range("I86").Value = 111113350
$ws.Range("J86").Value = 1635.4
$ws.Range("K86").Value = 111113350
$ws.Range("L86").Value = 1635.4
$ws.Range("M86").Value = -111112227
$ws.Range("N86").Value = -3881.4
$ws.Range("H89").Value = 71430590
$ws.Range("I89").Value = 111113350
$ws.Range("J89").Value = 1635.4
$ws.Range("K89").Value = 555566750
$ws.Range("L89").Value = 8177
$ws.Range("M89").Value = -555561134
$ws.Range("N89").Value = -19409
$ws.Range("H105").Value = 41667996
$ws.Range("J105").Value = 2450
$ws.Range("L105").Value = 2450
$ws.Range("N105").Value = -5944
$ws.Range("H113").Value = 10990477
$ws.Range("I113").Value = 12822060
$ws.Range("J113").Value = 980
$ws.Range("K113").Value = 12822060
$ws.Range("L113").Value = 980
$ws.Range("M113").Value = -12819890
$ws.Range("N113").Value = -5320
$ws.Range("H134").Value = 2488.4517
$ws.Range("I134").Value = 2488.4517
$ws.Range("K134").Value = 7465.355100000001
$ws.Range("M134").Value = -4930.355100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5375
$ws.Range("J62").Value = 6750
$ws.Range("L62").Value = 20250
$ws.Range("N62").Value = -21622
$ws.Range("H65").Value = 5375
$ws.Range("J65").Value = 6750
$ws.Range("L65").Value = 60750
$ws.Range("N65").Value = -67614
$ws.Range("H68").Value = 3339.75
$ws.Range("I68").Value = 4877.76
$ws.Range("J68").Value = 2099.4194
$ws.Range("K68").Value = 14633.28
$ws.Range("L68").Value = 6298.2582
$ws.Range("M68").Value = -13822.28
$ws.Range("N68").Value = -7920.2582
$ws.Range("H71").Value = 3339.75
$ws.Range("I71").Value = 4877.76
$ws.Range("J71").Value = 2099.4194
$ws.Range("K71").Value = 43899.84
$ws.Range("L71").Value = 18894.7746
$ws.Range("M71").Value = -39843.84
$ws.Range("N71").Value = -27006.7746
$ws.Range("H76").Value = 3666.6667
$ws.Range("J76").Value = 3666.6667
$ws.Range("L76").Value = 11000.0001
$ws.Range("N76").Value = -11766.0001
$ws.Range("H79").Value = 3666.6667
$ws.Range("J79").Value = 3666.6667
$ws.Range("L79").Value = 11000.0001
$ws.Range("N79").Value = -13652.0001
$ws.Range("H122").Value = 493.53125
$ws.Range("I122").Value = 529.7143
$ws.Range("J122").Value = 483.4
$ws.Range("K122").Value = 4767.428699999999
$ws.Range("L122").Value = 4350.599999999999
$ws.Range("M122").Value = -2317.428699999999
$ws.Range("N122").Value = -9250.599999999999
$ws.Range("H131").Value = 13751151
$ws.Range("I131").Value = 7143261.5
$ws.Range("J131").Value = 15152824
$ws.Range("K131").Value = 21429784.5
$ws.Range("L131").Value = 45458472
$ws.Range("M131").Value = -21424744.5
$ws.Range("N131").Value = -45468552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6901.263
$ws.Range("I126").Value = 6901.263
$ws.Range("K126").Value = 20703.789
$ws.Range("M126").Value = -18233.789

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7146127.5
$ws.Range("I122").Value = 7939830.5
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 23819491.5
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -23817041.5
$ws.Range("N122").Value = -13300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9209.362999999999
$ws.Range("I62").Value = 3360
$ws.Range("J62").Value = 14083.833
$ws.Range("K62").Value = 3360
$ws.Range("L62").Value = 14083.833
$ws.Range("M62").Value = -2736
$ws.Range("N62").Value = -15331.833
$ws.Range("H65").Value = 9209.362999999999
$ws.Range("I65").Value = 3360
$ws.Range("J65").Value = 14083.833
$ws.Range("K65").Value = 16800
$ws.Range("L65").Value = 70419.16500000001
$ws.Range("M65").Value = -13680
$ws.Range("N65").Value = -76659.16500000001
$ws.Range("H122").Value = 5477.625
$ws.Range("I122").Value = 6049.1113
$ws.Range("J122").Value = 4742.857
$ws.Range("K122").Value = 18147.3339
$ws.Range("L122").Value = 14228.571
$ws.Range("M122").Value = -15697.3339
$ws.Range("N122").Value = -19128.571
